{"js": "// Update the date line and the 25 multiplication problems to the new\n// \"output generated at c986bee\" values.\nconst replacements = [\n  [\"2025-01-07 Tuesday\", \"2025-01-08 Wednesday\"],\n  [\"67\\u00D733=\", \"78\\u00D751=\"],\n  [\"44\\u00D743=\", \"59\\u00D711=\"],\n  [\"21\\u00D784=\", \"22\\u00D748=\"],\n  [\"33\\u00D771=\", \"77\\u00D762=\"],\n  [\"48\\u00D794=\", \"78\\u00D738=\"],\n  [\"52\\u00D794=\", \"37\\u00D748=\"],\n  [\"75\\u00D750=\", \"85\\u00D750=\"],\n  [\"15\\u00D764=\", \"79\\u00D771=\"],\n  [\"72\\u00D768=\", \"16\\u00D786=\"],\n  [\"45\\u00D754=\", \"81\\u00D712=\"],\n  [\"54\\u00D724=\", \"69\\u00D773=\"],\n  [\"55\\u00D741=\", \"96\\u00D798=\"],\n  [\"38\\u00D765=\", \"64\\u00D795=\"],\n  [\"55\\u00D734=\", \"89\\u00D744=\"],\n  [\"88\\u00D788=\", \"36\\u00D783=\"],\n  [\"52\\u00D727=\", \"47\\u00D742=\"],\n  [\"65\\u00D774=\", \"29\\u00D736=\"],\n  [\"85\\u00D767=\", \"21\\u00D731=\"],\n  [\"48\\u00D725=\", \"22\\u00D765=\"],\n  [\"47\\u00D781=\", \"22\\u00D788=\"],\n  [\"35\\u00D725=\", \"12\\u00D787=\"],\n  [\"90\\u00D781=\", \"31\\u00D785=\"],\n  [\"30\\u00D717=\", \"73\\u00D756=\"],\n  [\"67\\u00D738=\", \"46\\u00D749=\"],\n  [\"86\\u00D793=\", \"44\\u00D789=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems to the new\n# \"output generated at c986bee\" values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-07 Tuesday\", \"2025-01-08 Wednesday\"),\n    @(\"67\u00d733=\", \"78\u00d751=\"),\n    @(\"44\u00d743=\", \"59\u00d711=\"),\n    @(\"21\u00d784=\", \"22\u00d748=\"),\n    @(\"33\u00d771=\", \"77\u00d762=\"),\n    @(\"48\u00d794=\", \"78\u00d738=\"),\n    @(\"52\u00d794=\", \"37\u00d748=\"),\n    @(\"75\u00d750=\", \"85\u00d750=\"),\n    @(\"15\u00d764=\", \"79\u00d771=\"),\n    @(\"72\u00d768=\", \"16\u00d786=\"),\n    @(\"45\u00d754=\", \"81\u00d712=\"),\n    @(\"54\u00d724=\", \"69\u00d773=\"),\n    @(\"55\u00d741=\", \"96\u00d798=\"),\n    @(\"38\u00d765=\", \"64\u00d795=\"),\n    @(\"55\u00d734=\", \"89\u00d744=\"),\n    @(\"88\u00d788=\", \"36\u00d783=\"),\n    @(\"52\u00d727=\", \"47\u00d742=\"),\n    @(\"65\u00d774=\", \"29\u00d736=\"),\n    @(\"85\u00d767=\", \"21\u00d731=\"),\n    @(\"48\u00d725=\", \"22\u00d765=\"),\n    @(\"47\u00d781=\", \"22\u00d788=\"),\n    @(\"35\u00d725=\", \"12\u00d787=\"),\n    @(\"90\u00d781=\", \"31\u00d785=\"),\n    @(\"30\u00d717=\", \"73\u00d756=\"),\n    @(\"67\u00d738=\", \"46\u00d749=\"),\n    @(\"86\u00d793=\", \"44\u00d789=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
